# Update column F (dSF) values for several rows to reflect the repulled
# data / recalculated mean, per the commit "repull data, push all data,
# mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    6  = 9
    7  = -8
    9  = -6
    10 = -2
    13 = -7
    14 = 6
    15 = -6
    17 = -2
    20 = -3
    25 = -4
    26 = 1
    28 = 4
    29 = -5
    30 = -4
    34 = 1
    36 = -1
    39 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
